# Agrícola del Norte S.A. de Arica - Plátano
# Weekly price update: a new observation is inserted as row 91 (pushing the
# existing rows 91-144 down to 92-145), for a new reporting date (44488)
# with updated min/max/average/kg prices, keeping the same variety/quality
# ("Sin especificar" / "Pintón"), volume (120) and origin (Ecuador) as the
# row that used to occupy row 91.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 91; existing rows 91..144 shift to 92..145.
$ws.Rows.Item(91).Insert()

# Populate the newly inserted row 91 with the new weekly observation.
$ws.Cells.Item(91, 1).Value = 1
$ws.Cells.Item(91, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(91, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(91, 4).Value = 44488
$ws.Cells.Item(91, 5).Value = 15
$ws.Cells.Item(91, 6).Value = "Fruta"
$ws.Cells.Item(91, 7).Value = 100108
$ws.Cells.Item(91, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(91, 9).Value = 100108006
$ws.Cells.Item(91, 10).Value = "Plátano"
$ws.Cells.Item(91, 11).Value = "Sin especificar"
$ws.Cells.Item(91, 12).Value = "Pintón"
$ws.Cells.Item(91, 13).Value = 120
$ws.Cells.Item(91, 14).Value = 23000
$ws.Cells.Item(91, 15).Value = 24000
$ws.Cells.Item(91, 16).Value = 23500
$ws.Cells.Item(91, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(91, 18).Value = "Ecuador"
$ws.Cells.Item(91, 19).Value = 1175
$ws.Cells.Item(91, 20).Value = 20
